$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Intrinsic analysis")

# Fill in the newly reported "Impact of static features" values
$ws.Range("C41").Value = 0.3664
$ws.Range("D41").Value = 0.514
$ws.Range("C42").Value = 0.2519
$ws.Range("D42").Value = 0.442
$ws.Range("C43").Value = 0.3584
$ws.Range("D43").Value = 0.467
$ws.Range("C44").Value = 0.2991
$ws.Range("D44").Value = 0.409

# Row height tweaks (several blocks switch from default to an explicit 15pt row height)
$ws.Range("A5:G5").RowHeight = 15
$ws.Range("A6:G6").RowHeight = 15
$ws.Range("A7:G7").RowHeight = 15
$ws.Range("A12:G12").RowHeight = 15
$ws.Range("A13:G13").RowHeight = 15
$ws.Range("A14:G14").RowHeight = 15

# Move the active selection/cursor
$ws.Range("C45").Select()
